$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date-format style (used by column A) down to the new row 53
$ws.Range("A52").Copy($ws.Range("A53"))

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 7.226520411029047
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 13.12477391005418
$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 7.424215059809214
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 14.69926045795804
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 4.268860212333636
$ws.Range("D4").Value = 2009
$ws.Range("E4").Value = 13.08276537368067
$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("C5").Value = -8.992252553594248
$ws.Range("D5").Value = 2010
$ws.Range("E5").Value = -19.76480035196674
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = -7.266312015249776
$ws.Range("D6").Value = 2010
$ws.Range("E6").Value = 12.31225042954258
$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 7.007132997505194
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 26.6762552377682
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 6.958243460951929
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 12.04357532583245
$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 8.86644397614711
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 5.982869268853186
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 9.469137444079934
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 10.6992064911972
$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 3.001306214623578
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 0.08911682035466217
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 3.358206407534947
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 4.390489499870109
$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = -2.904769335987201
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = -8.13908495190001
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.3081076735359067
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 3.50243535103556
$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 6.148460028297587
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 9.131012060398724
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 3.901355411819707
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 6.143002545701282
$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 5.171596082708629
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 5.986046065677453
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 5.331683351557981
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 4.555278923792572
$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 4.337699953939178
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 5.703761500281579
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 3.254758369308375
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 0.9515943257393689
$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 4.073887526082065
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 1.602827009554897
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 5.246209615995667
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 4.251116704684899
$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 7.340964210079881
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 8.197760099691198
$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 2.281540236993274
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = -4.308894244053674
$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 5.15263050486201
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 6.518301903862911
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 4.86255966374296
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 4.992093705734701
$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 3.497157880977619
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 2.970919592955878
$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 3.436682959168125
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 2.82953744009995
$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 2.727571990295941
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = -0.7915059299106297
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 2.764740011159428
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 0.3611963426346065
$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 2.93113923573054
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 5.185727774795068
$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -2.305533699949835
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -6.248031846400004
$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -7.578477024949737
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = 17.08516853885251
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = -7.260793671746435
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 2.387971016884638
$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = -0.8212867190516282
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = -3.88653222228782
$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 5.121832664816339
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 12.25634856540583
$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 4.379227219808146
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 3.504204983279191
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 4.097586525396268
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 3.347989317130651
$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 5.102173676573241
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = -5.324806173859886
$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 8.333485306093348
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 7.042341419899389
$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 7.397318165265498
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 1.541970075766486
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 7.824284864703746
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = -2.122471977790918
$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -0.9123477982208139
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = -5.598178600215808
$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = 0.08405665459807476
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 0.4484505192704713
$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = -0.3046246622258053
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 1.814369137627603
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = -1.24502235313334
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = -4.678511595261359
$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = -3.378144228902036
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = -0.366247000203368
$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = -2.055826728150212
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 3.376972582720295
$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = -2.567041707495976
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = -0.6241557099098238
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = -1.735114423676209
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 1.758477003221981
$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 1.830872485486124
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = -0.1271413384857256
$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 1.378024997308636
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 0.4446401485209472
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = 2.64031107104763
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 2.021190955273178
